$d = $word.ActiveDocument

# The document contains three "<id>...</id>" tag-text paragraphs, each
# previously split across three runs (the literal "<id>" in Courier New,
# the id value in the default font, and the literal "</id>" in Courier
# New again). Collapse each of them into a single run of Courier New
# text, e.g. "<id>p037v_1</id>", by doing a literal find/replace on the
# full visible text (which also merges the three runs into one, taking
# on the formatting of the first/matched run).
$d.Content.Find.Execute("<id>p037v_1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p037v_1</id>", 2) | Out-Null
$d.Content.Find.Execute("<id>p037v_2</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p037v_2</id>", 2) | Out-Null
$d.Content.Find.Execute("<id>p037v_3</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p037v_3</id>", 2) | Out-Null
